$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 57
$ws.Range("I11").Value = 57
$ws.Range("K11").Value = 57
$ws.Range("M11").Value = 83

$ws.Range("H98").Value = 3101.0908
$ws.Range("I98").Value = 1324.9231
$ws.Range("K98").Value = 1324.9231
$ws.Range("M98").Value = 173.0769

$ws.Range("H122").Value = 3101.0908
$ws.Range("I122").Value = 1324.9231
$ws.Range("K122").Value = 3974.7693
$ws.Range("M122").Value = -1524.7693

$ws.Range("H123").Value = 71670.3
$ws.Range("J123").Value = 71670.3
$ws.Range("L123").Value = 71670.3
$ws.Range("N123").Value = -81470.3

$ws.Range("H125").Value = 7152.8696
$ws.Range("J125").Value = 8673.857
$ws.Range("L125").Value = 78064.713
$ws.Range("N125").Value = -82984.713

$ws.Range("H128").Value = 87246.53999999999
$ws.Range("J128").Value = 87246.53999999999
$ws.Range("L128").Value = 87246.53999999999
$ws.Range("N128").Value = -97206.53999999999

$ws.Range("H134").Value = 65499.7
$ws.Range("J134").Value = 65499.7
$ws.Range("L134").Value = 65499.7
$ws.Range("N134").Value = -75639.7

$ws.Range("H137").Value = 3441.4893
$ws.Range("I137").Value = 1840.1714
$ws.Range("J137").Value = 8112
$ws.Range("K137").Value = 5520.5142
$ws.Range("L137").Value = 24336
$ws.Range("M137").Value = -2970.5142
$ws.Range("N137").Value = -29436

$ws.Range("H138").Value = 3782.0483
$ws.Range("I138").Value = 2145.077
$ws.Range("J138").Value = 4964.3057
$ws.Range("K138").Value = 6435.231000000001
$ws.Range("L138").Value = 14892.9171
$ws.Range("M138").Value = -1295.231000000001
$ws.Range("N138").Value = -25172.9171

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2567.276
$ws.Range("I32").Value = 2078.25
$ws.Range("K32").Value = 2078.25
$ws.Range("M32").Value = -1791.25

$ws.Range("H45").Value = 3300.3572
$ws.Range("I45").Value = 3423.125
$ws.Range("K45").Value = 3423.125
$ws.Range("M45").Value = -3046.125

$ws.Range("H61").Value = 2927.1428
$ws.Range("I61").Value = 2086.8235
$ws.Range("K61").Value = 2086.8235
$ws.Range("M61").Value = -1874.8235

$ws.Range("H63").Value = 3941
$ws.Range("I63").Value = 3941
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 3941
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -3255
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 3941
$ws.Range("I66").Value = 3941
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 19705
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -16273
$ws.Range("N66").ClearContents()

$ws.Range("H74").Value = 2777.923
$ws.Range("I74").Value = 1685.6666
$ws.Range("K74").Value = 1685.6666
$ws.Range("M74").Value = -811.6666

$ws.Range("H77").Value = 2777.923
$ws.Range("I77").Value = 1685.6666
$ws.Range("K77").Value = 8428.333000000001
$ws.Range("M77").Value = -4060.333000000001

$ws.Range("H102").Value = 2017.579
$ws.Range("I102").Value = 2046.3334
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 2046.3334
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = -424.3334
$ws.Range("N102").Value = -4744

$ws.Range("H122").Value = 3410.1538
$ws.Range("I122").Value = 2328.3
$ws.Range("K122").Value = 6984.900000000001
$ws.Range("M122").Value = -4534.900000000001

$ws.Range("H132").Value = 3490.45
$ws.Range("I132").Value = 3121.147
$ws.Range("K132").Value = 9363.440999999999
$ws.Range("M132").Value = -6833.440999999999

$ws.Range("H136").Value = 2927.1428
$ws.Range("I136").Value = 2086.8235
$ws.Range("K136").Value = 6260.470499999999
$ws.Range("M136").Value = -3710.470499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 95615.8
$ws.Range("J58").Value = 95615.8
$ws.Range("L58").Value = 95615.8
$ws.Range("N58").Value = -96203.8

$ws.Range("H59").Value = 81847.2
$ws.Range("J59").Value = 81847.2
$ws.Range("L59").Value = 81847.2
$ws.Range("N59").Value = -83541.2

$ws.Range("H74").Value = 89329.44500000001
$ws.Range("I74").Value = 89709
$ws.Range("J74").Value = 89282
$ws.Range("K74").Value = 89709
$ws.Range("L74").Value = 89282
$ws.Range("N74").Value = -91154
$ws.Range("M74").Value = -88773

$ws.Range("H77").Value = 89329.44500000001
$ws.Range("I77").Value = 89709
$ws.Range("J77").Value = 89282
$ws.Range("K77").Value = 269127
$ws.Range("L77").Value = 267846
$ws.Range("N77").Value = -277206
$ws.Range("M77").Value = -264447

$ws.Range("H82").Value = 98090.336
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 98090.336
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 98090.336
$ws.Range("N82").Value = -98856.336
$ws.Range("M82").ClearContents()

$ws.Range("H85").Value = 98090.336
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 98090.336
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 98090.336
$ws.Range("N85").Value = -100742.336
$ws.Range("M85").ClearContents()

$ws.Range("H86").Value = 2985
$ws.Range("I86").Value = 2676.1428
$ws.Range("K86").Value = 2676.1428
$ws.Range("M86").Value = -1553.1428

$ws.Range("H88").Value = 17749.75
$ws.Range("J88").Value = 17749.75
$ws.Range("L88").Value = 17749.75
$ws.Range("N88").Value = -18561.75

$ws.Range("H89").Value = 2985
$ws.Range("I89").Value = 2676.1428
$ws.Range("K89").Value = 13380.714
$ws.Range("M89").Value = -7764.714

$ws.Range("H91").Value = 17749.75
$ws.Range("J91").Value = 17749.75
$ws.Range("L91").Value = 17749.75
$ws.Range("N91").Value = -20557.75

$ws.Range("H99").Value = 1709.3448
$ws.Range("I99").Value = 1368.1538
$ws.Range("K99").Value = 1368.1538
$ws.Range("M99").Value = 129.8462

$ws.Range("H107").Value = 1026.6666
$ws.Range("I107").Value = 908
$ws.Range("K107").Value = 908
$ws.Range("M107").Value = 1012

$ws.Range("H134").Value = 3890.3103
$ws.Range("I134").Value = 2582.7273
$ws.Range("K134").Value = 7748.1819
$ws.Range("M134").Value = -5213.1819

$ws.Range("H138").Value = 49956
$ws.Range("J138").Value = 49956
$ws.Range("L138").Value = 49956
$ws.Range("N138").Value = -60236

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 600.0909
$ws.Range("I22").Value = 675.1667
$ws.Range("K22").Value = 675.1667
$ws.Range("M22").Value = -325.1667

$ws.Range("H31").Value = 3589.4707
$ws.Range("I31").Value = 1818.1428
$ws.Range("K31").Value = 1818.1428
$ws.Range("M31").Value = -1523.1428

$ws.Range("H34").Value = 3589.4707
$ws.Range("I34").Value = 1818.1428
$ws.Range("K34").Value = 1818.1428
$ws.Range("M34").Value = -1616.1428

$ws.Range("H94").Value = 2353.2222
$ws.Range("I94").Value = 680.6667
$ws.Range("K94").Value = 680.6667
$ws.Range("M94").Value = -229.6667

$ws.Range("H95").Value = 6906
$ws.Range("J95").Value = 6906
$ws.Range("L95").Value = 6906
$ws.Range("N95").Value = -12398

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1356.1428
$ws.Range("J98").Value = 961.4545000000001
$ws.Range("L98").Value = 2884.3635
$ws.Range("N98").Value = -5880.3635

$ws.Range("H112").Value = 1937
$ws.Range("I112").Value = 1937
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 5811
$ws.Range("L112").Value = 0
$ws.Range("M112").Value = -4703
$ws.Range("N112").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 41672816
$ws.Range("I70").Value = 5189.5713
$ws.Range("J70").Value = 58830076
$ws.Range("K70").Value = 5189.5713
$ws.Range("L70").Value = 58830076
$ws.Range("M70").Value = -4919.5713
$ws.Range("N70").Value = -58830616

$ws.Range("H73").Value = 41672816
$ws.Range("I73").Value = 5189.5713
$ws.Range("J73").Value = 58830076
$ws.Range("K73").Value = 5189.5713
$ws.Range("L73").Value = 58830076
$ws.Range("M73").Value = -4253.5713
$ws.Range("N73").Value = -58831948

$ws.Range("H80").Value = 1821423.6
$ws.Range("J80").Value = 2003562.4
$ws.Range("L80").Value = 2003562.4
$ws.Range("N80").Value = -2005558.4

$ws.Range("H83").Value = 1821423.6
$ws.Range("J83").Value = 2003562.4
$ws.Range("L83").Value = 10017812
$ws.Range("N83").Value = -10027796

$ws.Range("H102").Value = 5090
$ws.Range("I102").Value = 4985.8076
$ws.Range("K102").Value = 4985.8076
$ws.Range("M102").Value = -3363.8076

$ws.Range("H122").Value = 7262.8623
$ws.Range("I122").Value = 7229.4443
$ws.Range("J122").Value = 7317.5454
$ws.Range("K122").Value = 21688.3329
$ws.Range("L122").Value = 21952.6362
$ws.Range("M122").Value = -19238.3329
$ws.Range("N122").Value = -26852.6362

$ws.Range("H126").Value = 90911224
$ws.Range("I126").Value = 166668100
$ws.Range("J126").Value = 2979.6
$ws.Range("K126").Value = 500004300
$ws.Range("L126").Value = 8938.799999999999
$ws.Range("M126").Value = -500001830
$ws.Range("N126").Value = -13878.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 180240
$ws.Range("J6").Value = 180240
$ws.Range("L6").Value = 180240
$ws.Range("N6").Value = -180464

$ws.Range("H46").Value = 3517.3
$ws.Range("I46").Value = 3025.5
$ws.Range("K46").Value = 3025.5
$ws.Range("M46").Value = -2837.5

$ws.Range("H82").Value = 2802
$ws.Range("I82").Value = 3561
$ws.Range("K82").Value = 3561
$ws.Range("M82").Value = -3200

$ws.Range("H85").Value = 2802
$ws.Range("I85").Value = 3561
$ws.Range("K85").Value = 3561
$ws.Range("M85").Value = -2313

$ws.Range("H122").Value = 1619370.1
$ws.Range("I122").Value = 1431557.6
$ws.Range("K122").Value = 4294672.800000001
$ws.Range("M122").Value = -4292222.800000001

$ws.Range("H132").Value = 2967.4443
$ws.Range("I132").Value = 2061
$ws.Range("K132").Value = 6183
$ws.Range("M132").Value = -3653

$ws.Range("H136").Value = 4246.9565
$ws.Range("I136").Value = 4056.8948
$ws.Range("J136").Value = 5149.75
$ws.Range("K136").Value = 12170.6844
$ws.Range("L136").Value = 15449.25
$ws.Range("M136").Value = -9620.6844
$ws.Range("N136").Value = -20549.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

$ws.Range("H52").Value = 5000
$ws.Range("I52").Value = 5000
$ws.Range("K52").Value = 5000
$ws.Range("M52").Value = -4774

$ws.Range("H81").Value = 18840.572
$ws.Range("I81").Value = 2314.3333
$ws.Range("J81").Value = 31235.25
$ws.Range("K81").Value = 4628.6666
$ws.Range("L81").Value = 62470.5
$ws.Range("M81").Value = -3567.6666
$ws.Range("N81").Value = -64592.5

$ws.Range("H84").Value = 18840.572
$ws.Range("I84").Value = 2314.3333
$ws.Range("J84").Value = 31235.25
$ws.Range("K84").Value = 23143.333
$ws.Range("L84").Value = 312352.5
$ws.Range("M84").Value = -17839.333
$ws.Range("N84").Value = -322960.5

$ws.Range("H113").Value = 1018.619
$ws.Range("I113").Value = 1080.4615
$ws.Range("J113").Value = 918.125
$ws.Range("K113").Value = 3241.3845
$ws.Range("L113").Value = 2754.375
$ws.Range("M113").Value = -1071.3845
$ws.Range("N113").Value = -7094.375

$ws.Range("H126").Value = 2782.5676
$ws.Range("I126").Value = 2800.32
$ws.Range("K126").Value = 8400.960000000001
$ws.Range("M126").Value = -5930.960000000001

$ws.Range("H132").Value = 5566.6113
$ws.Range("I132").Value = 4519.9
$ws.Range("J132").Value = 6875
$ws.Range("K132").Value = 13559.7
$ws.Range("L132").Value = 20625
$ws.Range("M132").Value = -11029.7
$ws.Range("N132").Value = -25685

$ws.Range("H136").Value = 557099.5600000001
$ws.Range("I136").Value = 589399.5
$ws.Range("K136").Value = 1768198.5
$ws.Range("M136").Value = -1765648.5
